$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update accuracy values in column D
$ws.Range("D2").Value = 0.96
$ws.Range("D3").Value = 0.86
$ws.Range("D6").Value = 0.38
$ws.Range("D8").Value = 0.96

# Update the active selection to E8
$ws.Range("E8").Select()
